$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) for rows 2-51 store plain text values
# (e.g. "43.873.13", "  +2.49%  ") as inlineStr cells. Force the range to a
# text number format before writing so Excel doesn't reinterpret
# number-looking strings (like "2.25" or "0.440") as numeric/date values and
# strip significant formatting (leading/trailing zeros, padding spaces).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.873.13"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").Value = "2.230.68"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "268.38"
$ws.Range("E5").Value = "  +4.75%  "
$ws.Range("D6").Value = "88.50"
$ws.Range("E6").Value = "  +13.76%  "
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("D10").Value = "46.51"
$ws.Range("E10").Value = "  +9.89%  "
$ws.Range("D11").Value = "0.0927"
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("D12").Value = "7.57"
$ws.Range("E12").Value = "  +8.19%  "
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").Value = "2.561.10"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "14.80"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "2.208.30"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "43.836.59"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "70.11"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("D23").Value = "232.63"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").Value = "8.89"
$ws.Range("E24").Value = "  -4.28%  "
$ws.Range("D25").Value = "2.60"
$ws.Range("E25").Value = "  +18.41%  "
$ws.Range("D27").Value = "10.87"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("E28").Value = "  +6.00%  "
$ws.Range("D29").Value = "40.21"
$ws.Range("E29").Value = "  -5.32%  "
$ws.Range("D30").Value = "2.25"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").Value = "175.44"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").Value = "0.0906"
$ws.Range("E32").Value = "  +4.13%  "
$ws.Range("D33").Value = "20.64"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "5.43"
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  +18.21%  "
$ws.Range("E40").Value = "  -5.67%  "
$ws.Range("D41").Value = "65.62"
$ws.Range("E41").Value = "  +8.46%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").Value = "5.36"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").Value = "0.0995"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").Value = "101.30"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").Value = "8.32"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("E48").Value = "  +7.68%  "
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").Value = "0.440"
$ws.Range("E50").Value = "  -9.56%  "
$ws.Range("D51").Value = "1.51"
$ws.Range("E51").Value = "  +3.12%  "

# Restore the default (unstyled) look so only the cell contents change,
# matching the original workbook's styling for these data rows.
$dataRange.Style = "Normal"
